# Add function to calculate production rate:
# The unit label next to the calculated "grams per atom" value (C8) is
# simplified from "g/atoms" to just "g".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "g"

# Move/restore the active selection to C9 (just below the table).
[void]$ws.Range("C9").Select()
